$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.873.67'
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').Value = '1.630.41'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.26%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.517'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.15%  '

$ws.Range('E7').Value = '  -0.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.28'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.58%  '

$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('E10').Value = '  -0.98%  '

$ws.Range('E11').Value = '  -0.12%  '

$ws.Range('D12').Value = '1.861.52'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').Value = '1.625.17'
$ws.Range('E13').Value = '  -0.60%  '

$ws.Range('E14').Value = '  -1.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.556'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.02%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.23%  '

$ws.Range('D17').Value = '27.903.27'
$ws.Range('E17').Value = '  -0.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.10%  '

$ws.Range('E19').Value = '  -0.45%  '

$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('E21').Value = '  -0.20%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.56%  '

$ws.Range('E24').Value = '  -0.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.17%  '

$ws.Range('E29').Value = '  -0.19%  '

$ws.Range('E30').Value = '  +0.05%  '

$ws.Range('E31').Value = '  -0.18%  '

$ws.Range('E32').Value = '  +0.69%  '

$ws.Range('D33').Value = '1.417.80'
$ws.Range('E33').Value = '  +1.04%  '

$ws.Range('E34').Value = '  +0.98%  '

$ws.Range('E35').Value = '  +2.37%  '

$ws.Range('E36').Value = '  -3.77%  '

$ws.Range('E37').Value = '  -1.50%  '

$ws.Range('E38').Value = '  -0.79%  '

$ws.Range('E39').Value = '  -0.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.854'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.16%  '

$ws.Range('E43').Value = '  -0.95%  '

$ws.Range('E44').Value = '  -0.60%  '

$ws.Range('D45').Value = '1.770.91'
$ws.Range('E45').Value = '  -0.37%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.55%  '

$ws.Range('E48').Value = '  +1.22%  '

$ws.Range('E49').Value = '  -0.54%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.91%  '

$ws.Range('E51').Value = '  -0.28%  '
